$wb = $excel.ActiveWorkbook

# --- Sheet1: append 5 new TODO rows (A10:A14) ---
$ws1 = $wb.Worksheets.Item(1)

$s1 = "RefSource selection doesn" + [char]0xFD + " need to be on output page, put it to Source page as checkbox or something (only needed for RefSync mode - no need to be on outpt page that is for all)"
$s2 = "Make Source name edit focused control when new Source is created"
$s3 = "Custom naming patterns - as advanced option"
$s4 = "Add seconds to default pattern - when adding additional prhotos to previously done mix there maybe overwrites since the counter starts at zero on additional mix"
$s5 = "When Tool is executed go to output page to see the log"

$ws1.Range("A10").Value = $s1
$ws1.Range("A11").Value = $s2
$ws1.Range("A12").Value = $s3
$ws1.Range("A13").Value = $s4
$ws1.Range("A14").Value = $s5

# Move the active selection down to the new last row, as in the author's edit
[void]$ws1.Range("A14").Select()

# --- Sheet2: move the (still empty) sheet's selection to B27 ---
$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Range("B27").Select()

# Re-activate Sheet1 so it stays the tab that is selected when the file is reopened
[void]$ws1.Activate()
[void]$ws1.Range("A14").Select()
